$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (E1, F1)
$ws.Range("E1").Value = "Execution Time (ms)"
$ws.Range("F1").Value = "Memory Usage (B)"

# Copy the header style (bold, centered, bordered) from an existing header cell
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# Execution Time (ms) values
$ws.Range("E2").Value = 3.831899986835197
$ws.Range("E3").Value = 3.197999991243705
$ws.Range("E4").Value = 36.15679999347776
$ws.Range("E5").Value = 2.02829999034293

# Memory Usage (B) values
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 4096
$ws.Range("F5").Value = 0
